$d = $word.ActiveDocument

$d.Content.Find.Execute("622÷3=207, 1", $true, $false, $false, $false, $false, $true, 1, $false, "350÷4=87, 2", 2) | Out-Null
$d.Content.Find.Execute("773÷5=154, 3", $true, $false, $false, $false, $false, $true, 1, $false, "453÷2=226, 1", 2) | Out-Null
$d.Content.Find.Execute("464÷4=116, 0", $true, $false, $false, $false, $false, $true, 1, $false, "418÷9=46, 4", 2) | Out-Null
$d.Content.Find.Execute("757÷3=252, 1", $true, $false, $false, $false, $false, $true, 1, $false, "848÷7=121, 1", 2) | Out-Null
$d.Content.Find.Execute("202÷5=40, 2", $true, $false, $false, $false, $false, $true, 1, $false, "535÷4=133, 3", 2) | Out-Null
$d.Content.Find.Execute("842÷3=280, 2", $true, $false, $false, $false, $false, $true, 1, $false, "993÷2=496, 1", 2) | Out-Null
$d.Content.Find.Execute("180÷4=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "663÷2=331, 1", 2) | Out-Null
$d.Content.Find.Execute("125÷7=17, 6", $true, $false, $false, $false, $false, $true, 1, $false, "183÷8=22, 7", 2) | Out-Null
$d.Content.Find.Execute("687÷6=114, 3", $true, $false, $false, $false, $false, $true, 1, $false, "732÷6=122, 0", 2) | Out-Null
$d.Content.Find.Execute("202÷8=25, 2", $true, $false, $false, $false, $false, $true, 1, $false, "760÷7=108, 4", 2) | Out-Null
$d.Content.Find.Execute("423÷2=211, 1", $true, $false, $false, $false, $false, $true, 1, $false, "762÷2=381, 0", 2) | Out-Null
$d.Content.Find.Execute("142÷6=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "176÷3=58, 2", 2) | Out-Null
$d.Content.Find.Execute("344÷5=68, 4", $true, $false, $false, $false, $false, $true, 1, $false, "903÷4=225, 3", 2) | Out-Null
$d.Content.Find.Execute("816÷8=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "944÷6=157, 2", 2) | Out-Null
$d.Content.Find.Execute("215÷4=53, 3", $true, $false, $false, $false, $false, $true, 1, $false, "371÷4=92, 3", 2) | Out-Null
$d.Content.Find.Execute("955÷8=119, 3", $true, $false, $false, $false, $false, $true, 1, $false, "136÷8=17, 0", 2) | Out-Null
$d.Content.Find.Execute("792÷4=198, 0", $true, $false, $false, $false, $false, $true, 1, $false, "809÷8=101, 1", 2) | Out-Null
$d.Content.Find.Execute("937÷5=187, 2", $true, $false, $false, $false, $false, $true, 1, $false, "761÷2=380, 1", 2) | Out-Null
$d.Content.Find.Execute("727÷9=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "107÷5=21, 2", 2) | Out-Null
$d.Content.Find.Execute("626÷9=69, 5", $true, $false, $false, $false, $false, $true, 1, $false, "264÷4=66, 0", 2) | Out-Null
$d.Content.Find.Execute("771÷4=192, 3", $true, $false, $false, $false, $false, $true, 1, $false, "769÷2=384, 1", 2) | Out-Null
$d.Content.Find.Execute("666÷7=95, 1", $true, $false, $false, $false, $false, $true, 1, $false, "244÷8=30, 4", 2) | Out-Null
$d.Content.Find.Execute("150÷4=37, 2", $true, $false, $false, $false, $false, $true, 1, $false, "686÷3=228, 2", 2) | Out-Null
$d.Content.Find.Execute("272÷4=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "836÷5=167, 1", 2) | Out-Null
$d.Content.Find.Execute("507÷7=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "488÷3=162, 2", 2) | Out-Null
